$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain plain text (they were stored as
# inlineStr in the source workbook, e.g. "68.318.88"), matching the
# pre-existing formatting of the untouched D cells in this sheet.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.353.41"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.446.41"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "556.93"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "161.94"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.156"
$ws.Range("E9").Value = "  +7.06%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "4.81"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("B13").Value = "WrappedBTC"
$ws.Range("C13").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D13").Value = "68.288.69"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000169"
$ws.Range("E14").Value = "  +3.38%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "23.26"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "10.45"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "336.37"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "6.89"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "3.78"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "1.89"
$ws.Range("E20").Value = "  +5.13%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "66.66"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("B23").Value = "NEARProtocol"
$ws.Range("C23").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D23").Value = "3.69"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("B24").Value = "Aptos"
$ws.Range("C24").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D24").Value = "8.14"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0₃0815"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "7.19"
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "425.95"
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.14"
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.61"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "160.75"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("B32").Value = "WhiteBITCoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D32").Value = "18.98"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "17.79"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.105"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "0.296"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "4.36"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "1.46"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.06"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.35"
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "129.45"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0716"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.481"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.561"
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0916"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("B47").Value = "BitgetToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("B48").Value = "Optimism"
$ws.Range("C48").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D48").Value = "1.36"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "16.71"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "4.89"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0205"
$ws.Range("E51").Value = "  +6.76%  "
